$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B. This shifts the existing B/C/D columns
# (withholding_tax_id, tax, base) one position to the right, to C/D/E,
# while preserving their exact original column widths.
$ws.Columns.Item(2).Insert()

# --- Column widths ---
# Column A keeps its original width (untouched).
# New column B ("_requirements") and column C ("withholding_tax_id")
# get new custom widths.
$ws.Columns.Item(2).ColumnWidth = 14.046666666666667
$ws.Columns.Item(3).ColumnWidth = 18.746666666666666
# Columns D and E keep the exact widths inherited from the old C and D
# columns (4.07 and 5.46) thanks to the column insert above.

# --- Header row ---
$ws.Range("B1").Value = "_requirements"

# --- Row 2 (unchanged data, already correct after the column insert) ---
$ws.Range("A2").Value = "z0bug.wt_1040_1"
$ws.Range("C2").Value = "z0bug.wt_1040"
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 1

# --- Row 3 ---
$ws.Range("A3").Value = "z0bug.wt_1038_1"
$ws.Range("B3").Value = "G=='zero'"
$ws.Range("C3").Value = "z0bug.wt_1038"
$ws.Range("D3").Value = 23
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.5"

# --- Row 4 ---
$ws.Range("A4").Value = "z0bug.wt_1038_1"
$ws.Range("B4").Value = "G!='zero'"
$ws.Range("C4").Value = "z0bug.wt_1038"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "11.5"
$ws.Range("E4").Value = 1

# --- Row 5 ---
$ws.Range("A5").Value = "z0bug.wt_enasarco_1_1"
$ws.Range("B5").Value = "G=='zero'"
$ws.Range("C5").Value = "z0bug.wt_enasarco_1"
$ws.Range("D5").Value = 17
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.5"

# --- Row 6 ---
$ws.Range("A6").Value = "z0bug.wt_enasarco_1_1"
$ws.Range("B6").Value = "G!='zero'"
$ws.Range("C6").Value = "z0bug.wt_enasarco_1"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.5"
$ws.Range("E6").Value = 1

# --- Row 7 ---
$ws.Range("A7").Value = "z0bug.wt_1040-23A_1"
$ws.Range("C7").Value = "z0bug.wt_1040-23A"
$ws.Range("D7").Value = 23
$ws.Range("E7").Value = 1

# --- Selection ---
$ws.Range("D5").Select() | Out-Null
